$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from SCD0202 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-217" to "SCD0011-033"
$ws.Range("B2").Value = "SCD0011-033"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.592447916666666

# Update the active selection to B3 (matches the saved UI state in the diff)
$ws.Range("B3").Select()
